$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.103.33"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.850.73"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'0.6947"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.35%  "
$ws.Range("D6").Value = "'238.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "'0.9996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.07736"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.53%  "
$ws.Range("D9").Value = "'0.3035"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").Value = "'23.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").Value = "'0.08116"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'0.7269"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.830.38"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.208"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'89.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").Value = "29.107.41"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "'5.754"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.62%  "
$ws.Range("D18").Value = "'13.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "'0.000007742"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'236.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "2.093.75"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "'0.9995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'7.599"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").Value = "'8.980"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "'161.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "'0.1436"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.16%  "
$ws.Range("D28").Value = "'18.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'1.402"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("D31").Value = "'4.503"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").Value = "'1.489"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").Value = "'0.05239"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'1.187"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7022"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.67%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.024"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("D40").Value = "'2.675"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").Value = "'0.9239"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.33%  "
$ws.Range("D42").Value = "'6.003"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "1.082.55"
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("D44").Value = "'0.4268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").Value = "'70.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "'0.9994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'103.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'1.779"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").Value = "1.990.48"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'9.150"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("D51").Value = "'7.010"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.68%  "
